$d = $word.ActiveDocument

# 1) Body: "RA: 000110364890 - 1" -> "RA:  " (clear the RA number, keep the run as two spaces)
$d.Content.Find.Execute(" 000110364890 - 1 ", $true, $false, $false, $false, $false, $true, 1, $false, "  ", 2)

# 2) Body: "A QWR," -> "A TERE," (bold run right after "A ")
$d.Content.Find.Execute("QWR", $true, $false, $false, $false, $false, $true, 1, $false, "TERE", 2)

# 3) Header: "DIRETORIA DE ENSINO REGIAO QWER" -> "...TRE"
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute("QWER", $true, $false, $false, $false, $false, $true, 1, $false, "TRE", 2)
    $hdr.Range.Find.Execute("QWR", $true, $false, $false, $false, $false, $true, 1, $false, "TERE", 2)

    # 5x "Qwer" -> "Tre" in the address line
    for ($i = 0; $i -lt 5; $i++) {
        $hdr.Range.Find.Execute("Qwer", $true, $false, $false, $false, $false, $true, 1, $false, "Tre", 2)
    }

    # 3x "qwer" -> "tre" in CEP / Tel / Email lines
    for ($i = 0; $i -lt 3; $i++) {
        $hdr.Range.Find.Execute("qwer", $true, $false, $false, $false, $false, $true, 1, $false, "tre", 2)
    }
}
